# Apply "Natmi following Dr Hou advice" update.
# Sheet is regenerated with ECs added as a 4th sending cluster (rows 2-25),
# in addition to already being a target cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Il1a"
$row2[0,2] = "Il1rap"
$row2[0,3] = "ECs"
$row2[0,4] = [double]"1"
$row2[0,5] = [double]"0.3333333333333333"
$row2[0,6] = [double]"0.038687"
$row2[0,7] = [double]"0.116061"
$row2[0,8] = [double]"0.0003724978460720701"
$row2[0,9] = [double]"0.0003724978460720701"
$row2[0,10] = [double]"3"
$row2[0,11] = [double]"1"
$row2[0,12] = [double]"1.393808333333333"
$row2[0,13] = [double]"4.181425"
$row2[0,14] = [double]"0.01389847670528152"
$row2[0,15] = [double]"0.01389847670528152"
$row2[0,16] = [double]"0.05392226299166666"
$row2[0,17] = [double]"0.485300366925"
$row2[0,18] = [double]"5.177152636400209E-06"
$row2[0,19] = [double]"5.177152636400209E-06"
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object "object[,]" 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Il1a"
$row3[0,2] = "Il1rap"
$row3[0,3] = "FAPs"
$row3[0,4] = [double]"1"
$row3[0,5] = [double]"0.3333333333333333"
$row3[0,6] = [double]"0.038687"
$row3[0,7] = [double]"0.116061"
$row3[0,8] = [double]"0.0003724978460720701"
$row3[0,9] = [double]"0.0003724978460720701"
$row3[0,10] = [double]"3"
$row3[0,11] = [double]"1"
$row3[0,12] = [double]"9.765272666666666"
$row3[0,13] = [double]"29.295818"
$row3[0,14] = [double]"0.09737523548435453"
$row3[0,15] = [double]"0.09737523548435455"
$row3[0,16] = [double]"0.3777891036553333"
$row3[0,17] = [double]"3.400101932898"
$row3[0,18] = [double]"3.627206547868267E-05"
$row3[0,19] = [double]"3.627206547868267E-05"
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object "object[,]" 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Il1a"
$row4[0,2] = "Il1rap"
$row4[0,3] = "M1"
$row4[0,4] = [double]"1"
$row4[0,5] = [double]"0.3333333333333333"
$row4[0,6] = [double]"0.038687"
$row4[0,7] = [double]"0.116061"
$row4[0,8] = [double]"0.0003724978460720701"
$row4[0,9] = [double]"0.0003724978460720701"
$row4[0,10] = [double]"3"
$row4[0,11] = [double]"1"
$row4[0,12] = [double]"7.943370666666667"
$row4[0,13] = [double]"23.830112"
$row4[0,14] = [double]"0.07920798687439083"
$row4[0,15] = [double]"0.07920798687439085"
$row4[0,16] = [double]"0.3073051809813333"
$row4[0,17] = [double]"2.765746628832"
$row4[0,18] = [double]"2.950480450241538E-05"
$row4[0,19] = [double]"2.950480450241539E-05"
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object "object[,]" 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Il1a"
$row5[0,2] = "Il1rap"
$row5[0,3] = "M2"
$row5[0,4] = [double]"1"
$row5[0,5] = [double]"0.3333333333333333"
$row5[0,6] = [double]"0.038687"
$row5[0,7] = [double]"0.116061"
$row5[0,8] = [double]"0.0003724978460720701"
$row5[0,9] = [double]"0.0003724978460720701"
$row5[0,10] = [double]"3"
$row5[0,11] = [double]"1"
$row5[0,12] = [double]"4.317150333333333"
$row5[0,13] = [double]"12.951451"
$row5[0,14] = [double]"0.04304882666150776"
$row5[0,15] = [double]"0.04304882666150777"
$row5[0,16] = [double]"0.1670175949456666"
$row5[0,17] = [double]"1.503158354511"
$row5[0,18] = [double]"1.603559520734155E-05"
$row5[0,19] = [double]"1.603559520734155E-05"
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object "object[,]" 1,20
$row6[0,0] = "ECs"
$row6[0,1] = "Il1a"
$row6[0,2] = "Il1rap"
$row6[0,3] = "Neutro"
$row6[0,4] = [double]"1"
$row6[0,5] = [double]"0.3333333333333333"
$row6[0,6] = [double]"0.038687"
$row6[0,7] = [double]"0.116061"
$row6[0,8] = [double]"0.0003724978460720701"
$row6[0,9] = [double]"0.0003724978460720701"
$row6[0,10] = [double]"3"
$row6[0,11] = [double]"1"
$row6[0,12] = [double]"71.75082300000001"
$row6[0,13] = [double]"215.252469"
$row6[0,14] = [double]"0.7154693498390701"
$row6[0,15] = [double]"0.7154693498390702"
$row6[0,16] = [double]"2.775824089401"
$row6[0,17] = [double]"24.982416804609"
$row6[0,18] = [double]"0.000266510791745638"
$row6[0,19] = [double]"0.000266510791745638"
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object "object[,]" 1,20
$row7[0,0] = "ECs"
$row7[0,1] = "Il1a"
$row7[0,2] = "Il1rap"
$row7[0,3] = "sCs"
$row7[0,4] = [double]"1"
$row7[0,5] = [double]"0.3333333333333333"
$row7[0,6] = [double]"0.038687"
$row7[0,7] = [double]"0.116061"
$row7[0,8] = [double]"0.0003724978460720701"
$row7[0,9] = [double]"0.0003724978460720701"
$row7[0,10] = [double]"3"
$row7[0,11] = [double]"1"
$row7[0,12] = [double]"5.114546"
$row7[0,13] = [double]"15.343638"
$row7[0,14] = [double]"0.0510001244353952"
$row7[0,15] = [double]"0.05100012443539521"
$row7[0,16] = [double]"0.197866441102"
$row7[0,17] = [double]"1.780797969918"
$row7[0,18] = [double]"1.899743650159226E-05"
$row7[0,19] = [double]"1.899743650159226E-05"
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object "object[,]" 1,20
$row8[0,0] = "M1"
$row8[0,1] = "Il1a"
$row8[0,2] = "Il1rap"
$row8[0,3] = "ECs"
$row8[0,4] = [double]"3"
$row8[0,5] = [double]"1"
$row8[0,6] = [double]"6.837337666666667"
$row8[0,7] = [double]"20.512013"
$row8[0,8] = [double]"0.06583331748909883"
$row8[0,9] = [double]"0.06583331748909883"
$row8[0,10] = [double]"3"
$row8[0,11] = [double]"1"
$row8[0,12] = [double]"1.393808333333333"
$row8[0,13] = [double]"4.181425"
$row8[0,14] = [double]"0.01389847670528152"
$row8[0,15] = [double]"0.01389847670528152"
$row8[0,16] = [double]"9.529938217613889"
$row8[0,17] = [double]"85.769443958525"
$row8[0,18] = [double]"0.0009149828295536428"
$row8[0,19] = [double]"0.0009149828295536428"
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object "object[,]" 1,20
$row9[0,0] = "M1"
$row9[0,1] = "Il1a"
$row9[0,2] = "Il1rap"
$row9[0,3] = "FAPs"
$row9[0,4] = [double]"3"
$row9[0,5] = [double]"1"
$row9[0,6] = [double]"6.837337666666667"
$row9[0,7] = [double]"20.512013"
$row9[0,8] = [double]"0.06583331748909883"
$row9[0,9] = [double]"0.06583331748909883"
$row9[0,10] = [double]"3"
$row9[0,11] = [double]"1"
$row9[0,12] = [double]"9.765272666666666"
$row9[0,13] = [double]"29.295818"
$row9[0,14] = [double]"0.09737523548435453"
$row9[0,15] = [double]"0.09737523548435455"
$row9[0,16] = [double]"66.76846662907045"
$row9[0,17] = [double]"600.916199661634"
$row9[0,18] = [double]"0.006410534793217275"
$row9[0,19] = [double]"0.006410534793217276"
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object "object[,]" 1,20
$row10[0,0] = "M1"
$row10[0,1] = "Il1a"
$row10[0,2] = "Il1rap"
$row10[0,3] = "M1"
$row10[0,4] = [double]"3"
$row10[0,5] = [double]"1"
$row10[0,6] = [double]"6.837337666666667"
$row10[0,7] = [double]"20.512013"
$row10[0,8] = [double]"0.06583331748909883"
$row10[0,9] = [double]"0.06583331748909883"
$row10[0,10] = [double]"3"
$row10[0,11] = [double]"1"
$row10[0,12] = [double]"7.943370666666667"
$row10[0,13] = [double]"23.830112"
$row10[0,14] = [double]"0.07920798687439083"
$row10[0,15] = [double]"0.07920798687439085"
$row10[0,16] = [double]"54.31150745949511"
$row10[0,17] = [double]"488.803567135456"
$row10[0,18] = [double]"0.005214524547574145"
$row10[0,19] = [double]"0.005214524547574146"
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object "object[,]" 1,20
$row11[0,0] = "M1"
$row11[0,1] = "Il1a"
$row11[0,2] = "Il1rap"
$row11[0,3] = "M2"
$row11[0,4] = [double]"3"
$row11[0,5] = [double]"1"
$row11[0,6] = [double]"6.837337666666667"
$row11[0,7] = [double]"20.512013"
$row11[0,8] = [double]"0.06583331748909883"
$row11[0,9] = [double]"0.06583331748909883"
$row11[0,10] = [double]"3"
$row11[0,11] = [double]"1"
$row11[0,12] = [double]"4.317150333333333"
$row11[0,13] = [double]"12.951451"
$row11[0,14] = [double]"0.04304882666150776"
$row11[0,15] = [double]"0.04304882666150777"
$row11[0,16] = [double]"29.51781458676255"
$row11[0,17] = [double]"265.660331280863"
$row11[0,18] = [double]"0.002834047073140223"
$row11[0,19] = [double]"0.002834047073140224"
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object "object[,]" 1,20
$row12[0,0] = "M1"
$row12[0,1] = "Il1a"
$row12[0,2] = "Il1rap"
$row12[0,3] = "Neutro"
$row12[0,4] = [double]"3"
$row12[0,5] = [double]"1"
$row12[0,6] = [double]"6.837337666666667"
$row12[0,7] = [double]"20.512013"
$row12[0,8] = [double]"0.06583331748909883"
$row12[0,9] = [double]"0.06583331748909883"
$row12[0,10] = [double]"3"
$row12[0,11] = [double]"1"
$row12[0,12] = [double]"71.75082300000001"
$row12[0,13] = [double]"215.252469"
$row12[0,14] = [double]"0.7154693498390701"
$row12[0,15] = [double]"0.7154693498390702"
$row12[0,16] = [double]"490.584604712233"
$row12[0,17] = [double]"4415.261442410098"
$row12[0,18] = [double]"0.04710172086167462"
$row12[0,19] = [double]"0.04710172086167463"
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object "object[,]" 1,20
$row13[0,0] = "M1"
$row13[0,1] = "Il1a"
$row13[0,2] = "Il1rap"
$row13[0,3] = "sCs"
$row13[0,4] = [double]"3"
$row13[0,5] = [double]"1"
$row13[0,6] = [double]"6.837337666666667"
$row13[0,7] = [double]"20.512013"
$row13[0,8] = [double]"0.06583331748909883"
$row13[0,9] = [double]"0.06583331748909883"
$row13[0,10] = [double]"3"
$row13[0,11] = [double]"1"
$row13[0,12] = [double]"5.114546"
$row13[0,13] = [double]"15.343638"
$row13[0,14] = [double]"0.0510001244353952"
$row13[0,15] = [double]"0.05100012443539521"
$row13[0,16] = [double]"34.96987801369933"
$row13[0,17] = [double]"314.728902123294"
$row13[0,18] = [double]"0.00335750738393892"
$row13[0,19] = [double]"0.00335750738393892"
$ws.Range("A13:T13").Value = $row13

$row14 = New-Object "object[,]" 1,20
$row14[0,0] = "M2"
$row14[0,1] = "Il1a"
$row14[0,2] = "Il1rap"
$row14[0,3] = "ECs"
$row14[0,4] = [double]"3"
$row14[0,5] = [double]"1"
$row14[0,6] = [double]"5.801338666666666"
$row14[0,7] = [double]"17.404016"
$row14[0,8] = [double]"0.05585819933486567"
$row14[0,9] = [double]"0.05585819933486567"
$row14[0,10] = [double]"3"
$row14[0,11] = [double]"1"
$row14[0,12] = [double]"1.393808333333333"
$row14[0,13] = [double]"4.181425"
$row14[0,14] = [double]"0.01389847670528152"
$row14[0,15] = [double]"0.01389847670528152"
$row14[0,16] = [double]"8.085954178088889"
$row14[0,17] = [double]"72.77358760279999"
$row14[0,18] = [double]"0.0007763438822546024"
$row14[0,19] = [double]"0.0007763438822546024"
$ws.Range("A14:T14").Value = $row14

$row15 = New-Object "object[,]" 1,20
$row15[0,0] = "M2"
$row15[0,1] = "Il1a"
$row15[0,2] = "Il1rap"
$row15[0,3] = "FAPs"
$row15[0,4] = [double]"3"
$row15[0,5] = [double]"1"
$row15[0,6] = [double]"5.801338666666666"
$row15[0,7] = [double]"17.404016"
$row15[0,8] = [double]"0.05585819933486567"
$row15[0,9] = [double]"0.05585819933486567"
$row15[0,10] = [double]"3"
$row15[0,11] = [double]"1"
$row15[0,12] = [double]"9.765272666666666"
$row15[0,13] = [double]"29.295818"
$row15[0,14] = [double]"0.09737523548435453"
$row15[0,15] = [double]"0.09737523548435455"
$row15[0,16] = [double]"56.65165391167644"
$row15[0,17] = [double]"509.8648852050879"
$row15[0,18] = [double]"0.00543920531396456"
$row15[0,19] = [double]"0.005439205313964561"
$ws.Range("A15:T15").Value = $row15

$row16 = New-Object "object[,]" 1,20
$row16[0,0] = "M2"
$row16[0,1] = "Il1a"
$row16[0,2] = "Il1rap"
$row16[0,3] = "M1"
$row16[0,4] = [double]"3"
$row16[0,5] = [double]"1"
$row16[0,6] = [double]"5.801338666666666"
$row16[0,7] = [double]"17.404016"
$row16[0,8] = [double]"0.05585819933486567"
$row16[0,9] = [double]"0.05585819933486567"
$row16[0,10] = [double]"3"
$row16[0,11] = [double]"1"
$row16[0,12] = [double]"7.943370666666667"
$row16[0,13] = [double]"23.830112"
$row16[0,14] = [double]"0.07920798687439083"
$row16[0,15] = [double]"0.07920798687439085"
$row16[0,16] = [double]"46.08218339219911"
$row16[0,17] = [double]"414.739650529792"
$row16[0,18] = [double]"0.004424415519743146"
$row16[0,19] = [double]"0.004424415519743147"
$ws.Range("A16:T16").Value = $row16

$row17 = New-Object "object[,]" 1,20
$row17[0,0] = "M2"
$row17[0,1] = "Il1a"
$row17[0,2] = "Il1rap"
$row17[0,3] = "M2"
$row17[0,4] = [double]"3"
$row17[0,5] = [double]"1"
$row17[0,6] = [double]"5.801338666666666"
$row17[0,7] = [double]"17.404016"
$row17[0,8] = [double]"0.05585819933486567"
$row17[0,9] = [double]"0.05585819933486567"
$row17[0,10] = [double]"3"
$row17[0,11] = [double]"1"
$row17[0,12] = [double]"4.317150333333333"
$row17[0,13] = [double]"12.951451"
$row17[0,14] = [double]"0.04304882666150776"
$row17[0,15] = [double]"0.04304882666150777"
$row17[0,16] = [double]"25.04525115857955"
$row17[0,17] = [double]"225.407260427216"
$row17[0,18] = [double]"0.002404629940790581"
$row17[0,19] = [double]"0.002404629940790581"
$ws.Range("A17:T17").Value = $row17

$row18 = New-Object "object[,]" 1,20
$row18[0,0] = "M2"
$row18[0,1] = "Il1a"
$row18[0,2] = "Il1rap"
$row18[0,3] = "Neutro"
$row18[0,4] = [double]"3"
$row18[0,5] = [double]"1"
$row18[0,6] = [double]"5.801338666666666"
$row18[0,7] = [double]"17.404016"
$row18[0,8] = [double]"0.05585819933486567"
$row18[0,9] = [double]"0.05585819933486567"
$row18[0,10] = [double]"3"
$row18[0,11] = [double]"1"
$row18[0,12] = [double]"71.75082300000001"
$row18[0,13] = [double]"215.252469"
$row18[0,14] = [double]"0.7154693498390701"
$row18[0,15] = [double]"0.7154693498390702"
$row18[0,16] = [double]"416.2508238350561"
$row18[0,17] = [double]"3746.257414515504"
$row18[0,18] = [double]"0.03996482956129752"
$row18[0,19] = [double]"0.03996482956129752"
$ws.Range("A18:T18").Value = $row18

$row19 = New-Object "object[,]" 1,20
$row19[0,0] = "M2"
$row19[0,1] = "Il1a"
$row19[0,2] = "Il1rap"
$row19[0,3] = "sCs"
$row19[0,4] = [double]"3"
$row19[0,5] = [double]"1"
$row19[0,6] = [double]"5.801338666666666"
$row19[0,7] = [double]"17.404016"
$row19[0,8] = [double]"0.05585819933486567"
$row19[0,9] = [double]"0.05585819933486567"
$row19[0,10] = [double]"3"
$row19[0,11] = [double]"1"
$row19[0,12] = [double]"5.114546"
$row19[0,13] = [double]"15.343638"
$row19[0,14] = [double]"0.0510001244353952"
$row19[0,15] = [double]"0.05100012443539521"
$row19[0,16] = [double]"29.67121347224533"
$row19[0,17] = [double]"267.040921250208"
$row19[0,18] = [double]"0.002848775116815258"
$row19[0,19] = [double]"0.002848775116815259"
$ws.Range("A19:T19").Value = $row19

$row20 = New-Object "object[,]" 1,20
$row20[0,0] = "Neutro"
$row20[0,1] = "Il1a"
$row20[0,2] = "Il1rap"
$row20[0,3] = "ECs"
$row20[0,4] = [double]"3"
$row20[0,5] = [double]"1"
$row20[0,6] = [double]"91.18095533333333"
$row20[0,7] = [double]"273.542866"
$row20[0,8] = [double]"0.8779359853299634"
$row20[0,9] = [double]"0.8779359853299635"
$row20[0,10] = [double]"3"
$row20[0,11] = [double]"1"
$row20[0,12] = [double]"1.393808333333333"
$row20[0,13] = [double]"4.181425"
$row20[0,14] = [double]"0.01389847670528152"
$row20[0,15] = [double]"0.01389847670528152"
$row20[0,16] = [double]"127.0887753848944"
$row20[0,17] = [double]"1143.79897846405"
$row20[0,18] = [double]"0.01220197284083688"
$row20[0,19] = [double]"0.01220197284083688"
$ws.Range("A20:T20").Value = $row20

$row21 = New-Object "object[,]" 1,20
$row21[0,0] = "Neutro"
$row21[0,1] = "Il1a"
$row21[0,2] = "Il1rap"
$row21[0,3] = "FAPs"
$row21[0,4] = [double]"3"
$row21[0,5] = [double]"1"
$row21[0,6] = [double]"91.18095533333333"
$row21[0,7] = [double]"273.542866"
$row21[0,8] = [double]"0.8779359853299634"
$row21[0,9] = [double]"0.8779359853299635"
$row21[0,10] = [double]"3"
$row21[0,11] = [double]"1"
$row21[0,12] = [double]"9.765272666666666"
$row21[0,13] = [double]"29.295818"
$row21[0,14] = [double]"0.09737523548435453"
$row21[0,15] = [double]"0.09737523548435455"
$row21[0,16] = [double]"890.4068908371542"
$row21[0,17] = [double]"8013.662017534388"
$row21[0,18] = [double]"0.08548922331169402"
$row21[0,19] = [double]"0.08548922331169403"
$ws.Range("A21:T21").Value = $row21

$row22 = New-Object "object[,]" 1,20
$row22[0,0] = "Neutro"
$row22[0,1] = "Il1a"
$row22[0,2] = "Il1rap"
$row22[0,3] = "M1"
$row22[0,4] = [double]"3"
$row22[0,5] = [double]"1"
$row22[0,6] = [double]"91.18095533333333"
$row22[0,7] = [double]"273.542866"
$row22[0,8] = [double]"0.8779359853299634"
$row22[0,9] = [double]"0.8779359853299635"
$row22[0,10] = [double]"3"
$row22[0,11] = [double]"1"
$row22[0,12] = [double]"7.943370666666667"
$row22[0,13] = [double]"23.830112"
$row22[0,14] = [double]"0.07920798687439083"
$row22[0,15] = [double]"0.07920798687439085"
$row22[0,16] = [double]"724.2841259534435"
$row22[0,17] = [double]"6518.557133580992"
$row22[0,18] = [double]"0.06953954200257112"
$row22[0,19] = [double]"0.06953954200257115"
$ws.Range("A22:T22").Value = $row22

$row23 = New-Object "object[,]" 1,20
$row23[0,0] = "Neutro"
$row23[0,1] = "Il1a"
$row23[0,2] = "Il1rap"
$row23[0,3] = "M2"
$row23[0,4] = [double]"3"
$row23[0,5] = [double]"1"
$row23[0,6] = [double]"91.18095533333333"
$row23[0,7] = [double]"273.542866"
$row23[0,8] = [double]"0.8779359853299634"
$row23[0,9] = [double]"0.8779359853299635"
$row23[0,10] = [double]"3"
$row23[0,11] = [double]"1"
$row23[0,12] = [double]"4.317150333333333"
$row23[0,13] = [double]"12.951451"
$row23[0,14] = [double]"0.04304882666150776"
$row23[0,15] = [double]"0.04304882666150777"
$row23[0,16] = [double]"393.6418917109518"
$row23[0,17] = [double]"3542.777025398566"
$row23[0,18] = [double]"0.03779411405236962"
$row23[0,19] = [double]"0.03779411405236963"
$ws.Range("A23:T23").Value = $row23

$row24 = New-Object "object[,]" 1,20
$row24[0,0] = "Neutro"
$row24[0,1] = "Il1a"
$row24[0,2] = "Il1rap"
$row24[0,3] = "Neutro"
$row24[0,4] = [double]"3"
$row24[0,5] = [double]"1"
$row24[0,6] = [double]"91.18095533333333"
$row24[0,7] = [double]"273.542866"
$row24[0,8] = [double]"0.8779359853299634"
$row24[0,9] = [double]"0.8779359853299635"
$row24[0,10] = [double]"3"
$row24[0,11] = [double]"1"
$row24[0,12] = [double]"71.75082300000001"
$row24[0,13] = [double]"215.252469"
$row24[0,14] = [double]"0.7154693498390701"
$row24[0,15] = [double]"0.7154693498390702"
$row24[0,16] = [double]"6542.308587092907"
$row24[0,17] = [double]"58880.77728383616"
$row24[0,18] = [double]"0.6281362886243523"
$row24[0,19] = [double]"0.6281362886243524"
$ws.Range("A24:T24").Value = $row24

$row25 = New-Object "object[,]" 1,20
$row25[0,0] = "Neutro"
$row25[0,1] = "Il1a"
$row25[0,2] = "Il1rap"
$row25[0,3] = "sCs"
$row25[0,4] = [double]"3"
$row25[0,5] = [double]"1"
$row25[0,6] = [double]"91.18095533333333"
$row25[0,7] = [double]"273.542866"
$row25[0,8] = [double]"0.8779359853299634"
$row25[0,9] = [double]"0.8779359853299635"
$row25[0,10] = [double]"3"
$row25[0,11] = [double]"1"
$row25[0,12] = [double]"5.114546"
$row25[0,13] = [double]"15.343638"
$row25[0,14] = [double]"0.0510001244353952"
$row25[0,15] = [double]"0.05100012443539521"
$row25[0,16] = [double]"466.3491903762786"
$row25[0,17] = [double]"4197.142713386508"
$row25[0,18] = [double]"0.04477484449813943"
$row25[0,19] = [double]"0.04477484449813944"
$ws.Range("A25:T25").Value = $row25

